# Update "paises.xlsx" COVID-19 country/provincia snapshot.
#
# The source data feed was refreshed (new pull at 18:16 instead of 17:46),
# which: 
#   - bumped several countries' case/death counters in place, and
#   - caused a handful of countries (India vs. Peru, Serbia vs. Taiwan/
#     Bulgaria/Argentina/Croacia, Lituania vs. Camboya, Ghana vs. Cuba/
#     Jamaica) to change rank order in the (descending, by total cases)
#     table, shifting which country's name+figures sit on which row.
#
# Rather than rely on a sort command, we just write the final, known-good
# country name + figures for every row whose contents actually changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 18:16"

# Each entry: row, Pais, Casos totales, Nuevos casos, Casos activos,
#             Recuperados, Casos criticos, Muertes hoy, Muertes
$updates = @(
    @(5,   "Italia",         47021, 5986, 5129, 37860, 2655, 627, 4032),
    @(8,   "Alemania",       18794, 3474,  180, 18561,    2,   9,   53),
    @(9,   "Estados Unidos", 16489, 2700,  125, 16141,   64,  16,  223),
    @(12,  "Suiza",           5164,  942,   15,  5095,    0,  11,   54),
    @(17,  "Noruega",         1908,  118,    1,  1900,   27,   0,    7),
    @(25,  "Chequia",          833,  139,    4,   829,    6,   0,    0),
    @(28,  "Brasil",           654,   14,    2,   645,   18,   0,    7),
    @(51,  "India",            249,   55,   23,   221,    0,   1,    5),
    @(52,  "Peru",             234,    0,    1,   230,    7,   2,    3),
    @(65,  "Serbia",           135,   32,    2,   133,    4,   0,    0),
    @(66,  "Taiwan",           135,   27,   28,   105,    0,   1,    2),
    @(67,  "Bulgaria",         129,   22,    1,   125,    0,   0,    3),
    @(68,  "Argentina",        128,    0,    3,   122,    0,   0,    3),
    @(69,  "Croacia",          128,   18,    5,   122,    0,   0,    1),
    @(92,  "Lituania",          51,    3,    1,    50,    1,   0,    0),
    @(93,  "Camboya",           51,   14,    1,    50,    0,   0,    0),
    @(115, "Ghana",             16,    5,    0,    16,    0,   0,    0),
    @(116, "Cuba",              16,    5,    0,    15,    0,   0,    1),
    @(117, "Jamaica",           16,    1,    2,    13,    0,   0,    1)
)

foreach ($u in $updates) {
    $row = $u[0]
    $ws.Range("A$row").Value = $u[1]
    $ws.Range("B$row").Value = $u[2]
    $ws.Range("C$row").Value = $u[3]
    $ws.Range("D$row").Value = $u[4]
    $ws.Range("E$row").Value = $u[5]
    $ws.Range("F$row").Value = $u[6]
    $ws.Range("G$row").Value = $u[7]
    $ws.Range("H$row").Value = $u[8]
}
